$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" date column (C2:C7) from 2023-09-16 (45185) to 2023-10-05 (45204)
$ws.Range("C2:C7").Value = 45204
